$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-09 05:48:37"
$ws.Range("H2").Value = "92%"
$ws.Range("E3").Value = "2026-02-09 05:48:39"
$ws.Range("E4").Value = "2026-02-09 05:48:42"
$ws.Range("L4").Value = "7.6 km/h - 104º 5:07 TU"
$ws.Range("N4").Value = "2.7 °C 5:19 TU"
$ws.Range("O4").Value = "4.3 °C"
$ws.Range("E5").Value = "2026-02-09 05:48:44"
$ws.Range("H5").Value = "89%"
$ws.Range("L5").Value = "19.4 km/h - 90º 5:12 TU"
$ws.Range("O5").Value = "-5.4 °C"
$ws.Range("E6").Value = "2026-02-09 05:48:47"
$ws.Range("N6").Value = "4.6 °C 5:02 TU"
$ws.Range("O6").Value = "6.4 °C"
$ws.Range("E7").Value = "2026-02-09 05:48:49"
$ws.Range("H7").Value = "70%"
$ws.Range("E8").Value = "2026-02-09 05:48:52"
$ws.Range("E9").Value = "2026-02-09 05:48:54"
$ws.Range("H9").Value = "88%"
$ws.Range("O9").Value = "6.7 °C"
$ws.Range("E10").Value = "2026-02-09 05:48:57"
$ws.Range("H10").Value = "95%"
$ws.Range("K10").Value = "-0.1 MJ/m2"
$ws.Range("N10").Value = "2.7 °C 5:29 TU"
$ws.Range("O10").Value = "5.4 °C"
$ws.Range("E11").Value = "2026-02-09 05:48:59"
$ws.Range("O11").Value = "2.2 °C"
$ws.Range("E12").Value = "2026-02-09 05:49:01"
$ws.Range("N12").Value = "4.2 °C 5:29 TU"
$ws.Range("O12").Value = "7.3 °C"
$ws.Range("E13").Value = "2026-02-09 05:49:04"
$ws.Range("J13").Value = "1011.2 hPa"
$ws.Range("L13").Value = "9.0 km/h - 15º 5:19 TU"
$ws.Range("N13").Value = "-3.7 °C 5:12 TU"
$ws.Range("O13").Value = "-1.8 °C"
$ws.Range("E14").Value = "2026-02-09 05:49:06"
$ws.Range("N14").Value = "6.2 °C 5:08 TU"
$ws.Range("O14").Value = "7.3 °C"
$ws.Range("E15").Value = "2026-02-09 05:49:09"
$ws.Range("N15").Value = "2.8 °C 5:24 TU"
$ws.Range("O15").Value = "5.4 °C"
$ws.Range("E16").Value = "2026-02-09 05:49:11"
$ws.Range("L16").Value = "37.8 km/h - 212º 5:15 TU"
$ws.Range("N16").Value = "-6.2 °C 5:29 TU"
$ws.Range("E17").Value = "2026-02-09 05:49:13"
$ws.Range("O17").Value = "-0.1 °C"
$ws.Range("E18").Value = "2026-02-09 05:49:16"
$ws.Range("N18").Value = "3.5 °C 5:29 TU"
$ws.Range("O18").Value = "6.1 °C"
$ws.Range("E19").Value = "2026-02-09 05:49:18"
$ws.Range("N19").Value = "2.8 °C 5:01 TU"
$ws.Range("O19").Value = "3.2 °C"
$ws.Range("E20").Value = "2026-02-09 05:49:21"
$ws.Range("E21").Value = "2026-02-09 05:49:23"
$ws.Range("H21").Value = "94%"
$ws.Range("J21").Value = "1010.0 hPa"
$ws.Range("K21").Value = "-0.1 MJ/m2"
$ws.Range("O21").Value = "0.5 °C"
$ws.Range("E22").Value = "2026-02-09 05:49:26"
$ws.Range("L22").Value = "26.6 km/h - 281º 5:15 TU"
$ws.Range("E23").Value = "2026-02-09 05:49:28"
$ws.Range("E24").Value = "2026-02-09 05:49:30"
$ws.Range("H24").Value = "89%"
$ws.Range("J24").Value = "1009.4 hPa"
$ws.Range("E25").Value = "2026-02-09 05:49:33"
$ws.Range("E26").Value = "2026-02-09 05:49:35"
$ws.Range("K26").Value = "-0.1 MJ/m2"
$ws.Range("E27").Value = "2026-02-09 05:49:37"
$ws.Range("E28").Value = "2026-02-09 05:49:40"
$ws.Range("O28").Value = "3.5 °C"
$ws.Range("E29").Value = "2026-02-09 05:49:42"
$ws.Range("H29").Value = "96%"
$ws.Range("K29").Value = "-0.1 MJ/m2"
$ws.Range("N29").Value = "2.6 °C 5:29 TU"
$ws.Range("O29").Value = "5.3 °C"
$ws.Range("E30").Value = "2026-02-09 05:49:45"
$ws.Range("O30").Value = "6.5 °C"
$ws.Range("E31").Value = "2026-02-09 05:49:47"
$ws.Range("O31").Value = "8.8 °C"
$ws.Range("E32").Value = "2026-02-09 05:49:50"
$ws.Range("E33").Value = "2026-02-09 05:49:52"
$ws.Range("N33").Value = "-2.2 °C 5:11 TU"
$ws.Range("O33").Value = "-0.7 °C"
$ws.Range("E34").Value = "2026-02-09 05:49:55"
$ws.Range("N34").Value = "-5.2 °C 5:09 TU"
$ws.Range("O34").Value = "-3.3 °C"
$ws.Range("E35").Value = "2026-02-09 05:49:57"
$ws.Range("J35").Value = "1009.9 hPa"
$ws.Range("N35").Value = "3.2 °C 5:28 TU"
$ws.Range("E36").Value = "2026-02-09 05:49:59"
$ws.Range("H36").Value = "84%"
$ws.Range("N36").Value = "5.3 °C 5:00 TU"
$ws.Range("O36").Value = "8.3 °C"
$ws.Range("E37").Value = "2026-02-09 05:50:02"
$ws.Range("H37").Value = "91%"
$ws.Range("J37").Value = "1009.4 hPa"
$ws.Range("N37").Value = "1.2 °C 5:27 TU"
$ws.Range("O37").Value = "3.2 °C"
$ws.Range("E38").Value = "2026-02-09 05:50:05"
$ws.Range("N38").Value = "3.7 °C 5:28 TU"
$ws.Range("O38").Value = "5.9 °C"
$ws.Range("E39").Value = "2026-02-09 05:50:07"
$ws.Range("E40").Value = "2026-02-09 05:50:10"
$ws.Range("E41").Value = "2026-02-09 05:50:12"
$ws.Range("E42").Value = "2026-02-09 05:50:14"
$ws.Range("N42").Value = "4.2 °C 5:24 TU"
$ws.Range("O42").Value = "6.3 °C"
$ws.Range("E43").Value = "2026-02-09 05:50:17"
$ws.Range("N43").Value = "5.9 °C 5:00 TU"
$ws.Range("E44").Value = "2026-02-09 05:50:19"
$ws.Range("L44").Value = "38.9 km/h - 210º 5:11 TU"
$ws.Range("O44").Value = "-7.3 °C"
$ws.Range("E45").Value = "2026-02-09 05:50:22"
$ws.Range("J45").Value = "1009.8 hPa"
$ws.Range("K45").Value = "-0.1 MJ/m2"
$ws.Range("E46").Value = "2026-02-09 05:50:24"
$ws.Range("H46").Value = "82%"
$ws.Range("O46").Value = "6.1 °C"
